$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowCH($row, $vals) {
    $arr = New-Object "object[,]" 1,6
    for ($i = 0; $i -lt 6; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("C$row`:H$row").Value = $arr
}

Set-RowCH 2 @(8.708540916442871, 0.2552618980407715, 13.9090576171875, 0.1481355747003131, -0.5275314320756641, -1.595921289858394)
Set-RowCH 3 @(24.69548416137696, -32.37852478027344, -13.10801219940186, 0.8970806586181024, -2.253284974872135, -2.632813992924737)
Set-RowCH 4 @(9.214871406555176, -8.127243041992188, -7.553183078765869, 3.44800701316113, -1.11263772954487, -4.982394183493419)
Set-RowCH 5 @(-1.296450614929199, -25.33206176757812, 4.173199653625488, -0.06218727720976736, 0.8519693432053494, -3.113853278584423)
Set-RowCH 6 @(3.13936185836792, 1.159258842468261, 10.52872467041016, -5.216508371043581, -1.343549091778454, 3.495817017180781)
Set-RowCH 7 @(-5.014933586120605, 3.021981239318848, 2.716179370880127, -5.247863911833437, 3.879994811812087, 4.198587553663369)
Set-RowCH 8 @(12.84034442901611, -46.83102416992188, 6.994163036346436, -3.767358408548382, 6.815083000828341, 0.4525070901939874)
Set-RowCH 9 @(39.45174789428711, 18.51260566711425, -32.22710418701172, 1.195290473119145, -1.12244814540705, -3.771863303259405)
Set-RowCH 10 @(-25.77109146118164, -4.628503799438477, -10.99930953979492, 4.371577180492927, -3.279605823973685, -4.797551259944596)
Set-RowCH 11 @(-25.15040969848633, -52.22837066650391, 40.054931640625, 0.8983272033211789, 1.015922405332796, -5.503626271692188)
Set-RowCH 12 @(-0.4132165908813476, 5.698812484741211, 7.912927150726318, -3.768744276456209, 1.106249593939454, 3.618037942816487)
Set-RowCH 13 @(17.88103866577148, 1.754395723342896, 14.25908279418945, -5.589900377533165, 3.503953431913461, 6.05607688553742)
Set-RowCH 14 @(-0.3170166015625, -36.44598770141602, 3.271512508392334, -2.11142427696607, 9.510962799581424, -0.04490738020751496)
Set-RowCH 15 @(10.27998065948486, 5.400233268737793, -10.95361423492432, 0.3853269297416066, -2.830002902065919, -0.4395293002041412)
Set-RowCH 16 @(-23.53899955749512, -6.704813957214356, -13.04527854919434, 7.810094511321666, -5.127505102706861, -1.701906724125914)
Set-RowCH 17 @(17.90056037902832, -19.72356796264648, -5.512706756591797, 4.436884702500196, -6.79071483312473, -0.1620084203350239)
Set-RowCH 18 @(-16.73113822937012, 13.58313274383545, 5.493848323822021, -0.8537680622172118, -3.271142317050352, 3.793036780432373)
Set-RowCH 19 @(4.27336597442627, -17.7934684753418, 25.06283569335937, -7.225753661849682, 1.326318728362095, 8.623776784117876)
Set-RowCH 20 @(-30.59529876708984, -20.34370040893555, -32.5362548828125, 3.725271386937927, 3.283059017820457, -3.099974366226723)
Set-RowCH 21 @(4.83027458190918, 5.810257911682129, -21.42179489135743, 2.260886770268454, -3.733322253401921, 0.4708228362480379)

# Append new rows 22-31
$ws.Range("A22").Value = 2000
$ws.Range("B22").Value = "walkingToRunning"
Set-RowCH 22 @(-11.62696266174316, -21.11956596374512, -16.67885398864746, 10.02548202544597, 5.10485092133122, 0.5539246020516877)
$ws.Range("A23").Value = 2100
$ws.Range("B23").Value = "walkingToRunning"
Set-RowCH 23 @(-13.2412109375, 3.214882850646973, 22.37087440490723, -2.205630293691351, -1.242692031161299, 1.271664639886122)
$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "walkingToRunning"
Set-RowCH 24 @(-15.93152809143066, 19.48006248474121, 7.360580921173096, -2.934371616203759, 0.8867856582421496, -0.5782981265278011)
$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "walkingToRunning"
Set-RowCH 25 @(32.82965087890625, -70.24890899658203, 12.97732639312744, -3.666905725004196, -0.9611553521678435, -2.932546598124879)
$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "walkingToRunning"
Set-RowCH 26 @(-13.05622291564941, 5.201236724853516, -37.13468551635742, 0.07597643337871512, 15.43764277153591, -3.051599523783981)
$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "walkingToRunning"
Set-RowCH 27 @(-35.03109741210938, -2.42187762260437, -11.56947135925293, 0.9658543077439407, -2.044685548512497, -3.529194889268377)
$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "walkingToRunning"
Set-RowCH 28 @(-34.99357986450195, -38.73427581787109, 22.45726585388184, 4.105489831320192, 13.38752136726634, -6.598481021626331)
$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "walkingToRunning"
Set-RowCH 29 @(-5.411000728607178, 11.04445457458496, 9.496528625488279, -1.887021554078156, 0.322081098381812, 1.926029883130062)
$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "walkingToRunning"
Set-RowCH 30 @(0.176581859588623, 12.73018646240234, 1.936479568481445, -4.947177827046182, 3.337713633454262, 1.887203427509456)
$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "walkingToRunning"
Set-RowCH 31 @(9.547159194946287, -75.66817474365234, 32.35793304443359, -5.072363108240451, -0.134192296503722, -3.107403220930567)

Write-Host "Edit complete"
